# The sales-data table ("producto_datos_venta") currently lives in columns
# B:C. This edit shifts the whole table one column to the left, into A:B,
# by deleting column A (everything to the right of it, including the data
# in B:C, shifts left automatically) and then re-sizing the table's range
# to match its new location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing table before we touch anything.
$lo = $ws.ListObjects.Item(1)

# Deleting column A shifts every column (including the table's data in
# B:C) one place to the left -> data ends up in A:B, and the per-column
# width formatting that was on column B now correctly lands on column A.
$ws.Columns("A").Delete()

# The ListObject/table definition itself isn't auto-updated by the column
# delete, so resize it explicitly to its new location.
$lo.Resize($ws.Range("A1:B6"))

# Match the author's final selection in the sheet.
$ws.Range("B11").Select()
